# "prese misure a distanza fissa"
# Add a new worksheet "d=10cm cambio volt" after the existing sheets,
# populate it with the fixed-distance voltage-change measurements, and
# make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "d=10cm cambio volt"

# Header row.
$ws.Range("A1").Value = "V"
$ws.Range("B1").Value = "teta1"
$ws.Range("C1").Value = "teta2"
$ws.Range("D1").Value = "teta3"

# Measurement rows.
$data = @(
    @(2, 8, 7, 6),
    @(3, 14, 17, 17),
    @(4, 33, 34, 34),
    @(5, 49, 49, 48),
    @(6, 70, 72, 73)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Make the new sheet the active one (matches activeTab="2" in the workbook).
$ws.Activate()
